$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the latest job match (replaces the old "Senior Software Engineer" entry)
$ws.Range("A2").Value = "Sr Data Scientist (Remote)"
$ws.Range("B2").Value = "First American"
$ws.Range("C2").Value = "Santa Rosa, CA, US USA"
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = "Data Scientist, RAG, MLflow, Python, SQL, R, Scala, Optimization, A/B Testing"

# F2 ("Posted At") must stay a literal text date string, not get auto-converted
# to a date serial number by the smart-typed Value setter.
$fCell = $ws.Range("F2")
$fCell.NumberFormat = "@"
$fCell.Value = "2026-02-20"
$fCell.ClearFormats()

$ws.Range("G2").Value = "https://www.indeed.com/viewjob?jk=7152380fec2fea5b"

# Remove the now-stale extra job matches (rows 3-5), shrinking the used range to A1:G2
$ws.Range("A3:G5").EntireRow.Delete()
